# "Solver may change input data" bug fix
# - Re-populate sheet "Лист1" (Sheet1) with the corrected solver iteration data
#   (a new column A with the iteration index 1..22 was introduced, and the
#   per-method value columns B/D/E/F were recomputed so the solvers no longer
#   mutate the shared input data).
# - Update the chart's category axis (tick label skip) and value axis
#   (cross-between) presentation to match the corrected, denser category axis.
# - Restore the worksheet selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Rewrite the data table (columns A-F, rows 2-23)
# ---------------------------------------------------------------------------

$ws.Range("A2").Value = 1.0
$ws.Range("B2").Value = -0.35046515366367464
$ws.Range("D2").Value = 0.17609125905568124
$ws.Range("E2").Value = 0.17609125905568124
$ws.Range("F2").Value = 0.17609125905568124

$ws.Range("A3").Value = 2.0
$ws.Range("B3").Value = -5.0
$ws.Range("D3").Value = -0.12493873660829995
$ws.Range("E3").Value = -0.42596873227228116
$ws.Range("F3").Value = -0.42596873227228116

$ws.Range("A4").Value = 3.0
$ws.Range("B4").Value = -5.0
$ws.Range("D4").Value = -0.42596873227228116
$ws.Range("E4").Value = -1.0280287236002434
$ws.Range("F4").Value = -1.0280287236002434

$ws.Range("A5").Value = 4.0
$ws.Range("B5").Value = -5.0
$ws.Range("D5").Value = -0.7269987279362623
$ws.Range("E5").Value = -1.630088714928206
$ws.Range("F5").Value = -1.630088714928206

$ws.Range("A6").Value = 5.0
$ws.Range("B6").Value = -5.0
$ws.Range("D6").Value = -1.0280287236002434
$ws.Range("E6").Value = -2.2321487062561682
$ws.Range("F6").Value = -2.2321487062561682

$ws.Range("A7").Value = 6.0
$ws.Range("D7").Value = -1.3290587192642247
$ws.Range("E7").Value = -2.8342086975841307
$ws.Range("F7").Value = -2.8342086975841307

$ws.Range("A8").Value = 7.0
$ws.Range("D8").Value = -1.630088714928206
$ws.Range("E8").Value = -3.4362686889120932
$ws.Range("F8").Value = -3.4362686889120932

$ws.Range("A9").Value = 8.0
$ws.Range("D9").Value = -1.9311187105921872
$ws.Range("E9").Value = -4.038328680240055
$ws.Range("F9").Value = -4.038328680240055

$ws.Range("A10").Value = 9.0
$ws.Range("D10").Value = -2.2321487062561682
$ws.Range("E10").Value = -4.640388671568018
$ws.Range("F10").Value = -4.640388671568018

$ws.Range("A11").Value = 10.0
$ws.Range("D11").Value = -2.5331787019201495
$ws.Range("E11").Value = -5.0
$ws.Range("F11").Value = -5.0

$ws.Range("A12").Value = 11.0
$ws.Range("D12").Value = -2.8342086975841307
$ws.Range("E12").Value = -5.0
$ws.Range("F12").Value = -5.0

$ws.Range("A13").Value = 12.0
$ws.Range("D13").Value = -3.135238693248112

$ws.Range("A14").Value = 13.0
$ws.Range("D14").Value = -3.4362686889120932

$ws.Range("A15").Value = 14.0
$ws.Range("D15").Value = -3.7372986845760745

$ws.Range("A16").Value = 15.0
$ws.Range("D16").Value = -4.038328680240055

$ws.Range("A17").Value = 16.0
$ws.Range("D17").Value = -4.339358675904037

$ws.Range("A18").Value = 17.0
$ws.Range("D18").Value = -4.640388671568018

$ws.Range("A19").Value = 18.0
$ws.Range("D19").Value = -4.941418667231999

$ws.Range("A20").Value = 19.0
$ws.Range("D20").Value = -5.0

$ws.Range("A21").Value = 20.0
$ws.Range("D21").Value = -5.0

$ws.Range("A22").Value = 21.0
$ws.Range("D22").Value = -5.0

$ws.Range("A23").Value = 22.0
$ws.Range("D23").Value = -5.0

# The tail of column F (rows 13-25) no longer has data - the corrected
# SeidelRelaxation series now converges and stops at row 12 (it shares
# column F with Seidel through row 12, see above).
$ws.Range("F13:F25").ClearContents()

# ---------------------------------------------------------------------------
# 2. Chart axis presentation: thin out the category labels and switch the
#    value axis to cross between categories at the midpoint.
# ---------------------------------------------------------------------------

$chart = $ws.ChartObjects(1).Chart
$catAx = $chart.Axes(1)
$valAx = $chart.Axes(2)

$catAx.TickLabelSpacing = 5
$valAx.AxisBetweenCategories = "midCat"

# ---------------------------------------------------------------------------
# 3. Restore the worksheet selection.
# ---------------------------------------------------------------------------

$ws.Range("A1:B4").Select()
